$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# "First page" footer (footer1.xml, docPr id="3") - Pearson logo: image1.png -> image2.png
$ftrFirst = $sec.Footers.Item(2)
$ftrFirst.Range.InlineShapes.Item(1).Name = "image2.png"

# "Default" footer (footer2.xml, docPr id="2") - Pearson logo: image1.png -> image2.png
$ftrDefault = $sec.Footers.Item(1)
$ftrDefault.Range.InlineShapes.Item(1).Name = "image2.png"

# "First page" header (header1.xml, docPr id="1") - BTec logo: image2.jpg -> image1.jpg
$hdrFirst = $sec.Headers.Item(2)
$hdrFirst.Range.InlineShapes.Item(1).Name = "image1.jpg"
